$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Foglio1")
$ws2 = $wb.Worksheets.Item("difficolta")

# Workbook-level: update the last-saved folder path shown in the file dialog history.
$wb.Path = "C:\APP - Copia\"

# New color options appended to the "difficolta" lookup sheet (column G).
$ws2.Range("G7").Value = "blue"
$ws2.Range("G8").Value = "white"
$ws2.Range("G9").Value = "yellow"

# Re-color some trail segments on the main sheet.
$ws1.Range("L2:L5").Value = "yellow"
$ws1.Range("L12:L15").Value = "orange"

# Restore selection on the lookup sheet, then leave the main sheet active/selected
# (matches the saved view state of the edited workbook).
[void]$ws2.Range("G16").Select()
[void]$ws1.Activate()
[void]$ws1.Range("L15").Select()
